# Delete the post row for "「ちょっと静けさが必要な時もある」" (row 406),
# shifting all subsequent rows up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(406).Delete()
